$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Martin Guptill / Temba Bavuma(C))
$ws.Range("B2").Value = 21
$ws.Range("C2").Value = 7
$ws.Range("D2").Value = 'Caught'
$ws.Range("E2").Value = ' Anrich Nortje'
$ws.Range("K2").Value = 51
$ws.Range("L2").Value = 19
$ws.Range("M2").Value = 'Bowled'
$ws.Range("N2").Value = ' Mitchell Santner'

# Row 3 (Daryl Mitchell / Quinton de Kock)
$ws.Range("B3").Value = 7
$ws.Range("C3").Value = 4
$ws.Range("D3").Value = 'LBW'
$ws.Range("E3").Value = ' Keshav Maharaj'
$ws.Range("K3").Value = 15
$ws.Range("L3").Value = 6
$ws.Range("N3").Value = ' Tim Southee'

# Row 4 (Kane Williamson(C) / Rassie Va der Dussen)
$ws.Range("B4").Value = 4
$ws.Range("C4").Value = 4
$ws.Range("D4").Value = 'Bowled'
$ws.Range("E4").Value = ' Keshav Maharaj'
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 'Bowled'
$ws.Range("N4").Value = ' Ish Sodhi'

# Row 5 (Devon Conway / Aiden Markram)
$ws.Range("B5").Value = 49
$ws.Range("C5").Value = 24
$ws.Range("D5").Value = 'Bowled'
$ws.Range("E5").Value = ' Kagiso Rabada'
$ws.Range("K5").Value = 5
$ws.Range("M5").Value = 'Caught'
$ws.Range("N5").Value = ' Ish Sodhi'

# Row 6 (Glenn Phillips / David Miller)
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 'LBW'
$ws.Range("E6").Value = ' Keshav Maharaj'
$ws.Range("K6").Value = 38
$ws.Range("L6").Value = 11
$ws.Range("M6").Value = 'Bowled'
$ws.Range("N6").Value = ' Trent Boult'

# Row 7 (James Neesham / Reeza Hendricks)
$ws.Range("B7").Value = 10
$ws.Range("C7").Value = 7
$ws.Range("E7").Value = ' Dwaine Pretorius'
$ws.Range("K7").Value = 10
$ws.Range("L7").Value = 3
$ws.Range("M7").Value = 'Caught'

# Row 8 (Mitchell Santner / Dwaine Pretorius)
$ws.Range("B8").Value = 18
$ws.Range("C8").Value = 6
$ws.Range("E8").Value = ' Tabraiz Shamsi'
$ws.Range("K8").Value = 9
$ws.Range("L8").Value = 4
$ws.Range("N8").Value = ' Ish Sodhi'

# Row 9 (Adam Milne / Kagiso Rabada)
$ws.Range("B9").Value = 12
$ws.Range("C9").Value = 4
$ws.Range("D9").Value = 'LBW'
$ws.Range("E9").Value = ' Anrich Nortje'
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 1
$ws.Range("N9").Value = ' Ish Sodhi'

# Row 10 (Ish Sodhi / Keshav Maharaj)
$ws.Range("B10").Value = 20
$ws.Range("C10").Value = 8
$ws.Range("D10").Value = 'Caught'
$ws.Range("E10").Value = ' Kagiso Rabada'
$ws.Range("K10").Value = 31
$ws.Range("L10").Value = 11
$ws.Range("M10").Value = 'LBW'
$ws.Range("N10").Value = ' Mitchell Santner'

# Row 11 (Tim Southee / Anrich Nortje)
$ws.Range("B11").Value = 12
$ws.Range("C11").Value = 5
$ws.Range("E11").Value = ' Tabraiz Shamsi'
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 3
$ws.Range("N11").Value = ' Trent Boult'

# Row 12 (Trent Boult / Tabraiz Shamsi)
$ws.Range("B12").Value = 40
$ws.Range("C12").Value = 10
$ws.Range("D12").Value = 'NOT OUT'
$ws.Range("E12").Value = ' '
$ws.Range("K12").Value = 8
$ws.Range("L12").Value = 3
$ws.Range("M12").Value = 'NOT OUT'
$ws.Range("N12").Value = ' '

# Row 16 (innings totals). C16/L16 hold overs like "13.2" which Excel would
# otherwise auto-convert to a number, so a leading quote forces text entry
# (same as a user typing '13.2 into the cell) while keeping General format.
$ws.Range("A16").Value = 193
$ws.Range("C16").Value = '''13.2'
$ws.Range("D16").Value = 80
$ws.Range("J16").Value = 169
$ws.Range("L16").Value = '''10.4'
$ws.Range("M16").Value = 64

# Row 21 (bowling figures)
$ws.Range("A21").Value = 'Dwaine Pretorius'
$ws.Range("C21").Value = 40
$ws.Range("D21").Value = 1
$ws.Range("J21").Value = 'Tim Southee'
$ws.Range("K21").Value = 12
$ws.Range("L21").Value = 30
$ws.Range("M21").Value = 1

# Row 22
$ws.Range("A22").Value = 'Tabraiz Shamsi'
$ws.Range("B22").Value = 18
$ws.Range("C22").Value = 40
$ws.Range("D22").Value = 2
$ws.Range("J22").Value = 'Ish Sodhi'
$ws.Range("L22").Value = 30

# Row 23
$ws.Range("A23").Value = 'Anrich Nortje'
$ws.Range("B23").Value = 18
$ws.Range("C23").Value = 54
$ws.Range("J23").Value = 'Adam Milne'
$ws.Range("L23").Value = 31
$ws.Range("M23").Value = 0

# Row 24
$ws.Range("A24").Value = 'Keshav Maharaj'
$ws.Range("B24").Value = 18
$ws.Range("C24").Value = 30
$ws.Range("D24").Value = 3
$ws.Range("J24").Value = 'Mitchell Santner'
$ws.Range("L24").Value = 32
$ws.Range("M24").Value = 2

# Row 25
$ws.Range("A25").Value = 'Kagiso Rabada'
$ws.Range("B25").Value = 14
$ws.Range("C25").Value = 29
$ws.Range("D25").Value = 2
$ws.Range("J25").Value = 'Trent Boult'
$ws.Range("K25").Value = 16
$ws.Range("L25").Value = 46
$ws.Range("M25").Value = 3
